$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert three new rows right after the current last row of the table (row 41),
# this mirrors the formatting of row 41 down into rows 42-44.
$ws.Rows("42:44").Insert()

# Expand ("resize") the table (ListObject) so it covers the newly inserted rows.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:T44"))

# Row 42: "General case study with files"
$ws.Range("A42").Value = "Case"
$ws.Range("B42").Value = 144018
$ws.Range("C42").Value = "General case study with files"
$ws.Range("D42").Value = "Membership testing"
$ws.Range("F42").Value = "Open government"
$ws.Range("G42").Value = "No"
$ws.Range("H42").Value = "Yes"
$ws.Range("P42").Value = "user9351@example.com"

# Row 43: "Guideline with files"
$ws.Range("A43").Value = "Case"
$ws.Range("B43").Value = 135110
$ws.Range("C43").Value = "Guideline with files"
$ws.Range("D43").Value = "Membership testing"
$ws.Range("F43").Value = "Open government"
$ws.Range("G43").Value = "No"
$ws.Range("H43").Value = "Yes"
$ws.Range("P43").Value = "user9351@example.com"

# Row 44: "Open source case study with files"
$ws.Range("A44").Value = "Case"
$ws.Range("B44").Value = 135160
$ws.Range("C44").Value = "Open source case study with files"
$ws.Range("D44").Value = "Membership testing"
$ws.Range("F44").Value = "Open government"
$ws.Range("G44").Value = "No"
$ws.Range("H44").Value = "Yes"
$ws.Range("P44").Value = "user9351@example.com"

# Update the hidden _FilterDatabase defined name so it covers the expanded table range.
$fdb = $wb.Names.Item("1. Content items!_FilterDatabase")
$fdb.RefersTo = "='1. Content items'!`$A`$1:`$T`$44"
